# Weekly CompStat report update: new crime data collected.
# Updates header (volume/week) and the crime-stat table (rows 15-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump volume "Number" and advance the reporting week ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# Row 15
$ws.Range("L15").Value = 80

# Row 16
$ws.Range("D16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -71.428571428571
$ws.Range("I16").Value = 149
$ws.Range("J16").Value = 143
$ws.Range("K16").Value = 4.195804195804
$ws.Range("L16").Value = 10.370370370370
$ws.Range("M16").Value = -21.164021164021
$ws.Range("N16").Value = -82.951945080091

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("C14").Copy($ws.Range("D17"))
$ws.Range("E14").Copy($ws.Range("E17"))
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 183
$ws.Range("J17").Value = 183
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 17.307692307692
$ws.Range("M17").Value = 79.411764705882
$ws.Range("N17").Value = -60.043668122270

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 180
$ws.Range("J18").Value = 133
$ws.Range("K18").Value = 35.338345864661
$ws.Range("L18").Value = -10
$ws.Range("M18").Value = 44
$ws.Range("N18").Value = -82.922201138519

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 2.702702702702
$ws.Range("I19").Value = 518
$ws.Range("J19").Value = 470
$ws.Range("K19").Value = 10.212765957446
$ws.Range("L19").Value = 51.020408163265
$ws.Range("M19").Value = -2.448210922787
$ws.Range("N19").Value = -53.249097472924

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 125
$ws.Range("I20").Value = 88
$ws.Range("J20").Value = 88
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 14.285714285714
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -91.020408163265

# Row 21
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -9.638554216867
$ws.Range("I21").Value = 1128
$ws.Range("J21").Value = 1032
$ws.Range("K21").Value = 9.302325581395
$ws.Range("L21").Value = 22.077922077922
$ws.Range("M21").Value = 12.350597609561
$ws.Range("N21").Value = -75.110326566637

# Row 22
$ws.Range("D16").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 33
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 17.857142857142
$ws.Range("L22").Value = 37.5
$ws.Range("M22").Value = 10

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 121
$ws.Range("J23").Value = 108
$ws.Range("K23").Value = 12.037037037037
$ws.Range("L23").Value = 27.368421052631
$ws.Range("M23").Value = 61.333333333333

# Row 24
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -19.444444444444
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 164
$ws.Range("H24").Value = -26.829268292682
$ws.Range("I24").Value = 1825
$ws.Range("J24").Value = 1296
$ws.Range("K24").Value = 40.817901234567
$ws.Range("L24").Value = 31.768953068592
$ws.Range("M24").Value = 72.985781990521

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = -12.5
$ws.Range("I25").Value = 301
$ws.Range("J25").Value = 273
$ws.Range("K25").Value = 10.256410256410
$ws.Range("L25").Value = 20.4
$ws.Range("M25").Value = -6.230529595015

# Row 26
$ws.Range("L26").Value = 50

# Row 27
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 53
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = 1.923076923076
$ws.Range("L27").Value = 35.897435897435

# Row 28
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("N28").Value = -90.909090909090

# Row 29
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("N29").Value = -92.307692307692
